# Update "想去人数" (people-interested count) figures to the latest
# scraped values across the "展览", "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14104
$ws1.Range("F4").Value = 550
$ws1.Range("F7").Value = 1049
$ws1.Range("F8").Value = 13934
$ws1.Range("F9").Value = 14985
$ws1.Range("F11").Value = 15
$ws1.Range("F19").Value = 24
$ws1.Range("F20").Value = 63
$ws1.Range("F22").Value = 1157
$ws1.Range("F25").Value = 5826
$ws1.Range("F26").Value = 948
$ws1.Range("F27").Value = 1067
$ws1.Range("F28").Value = 5448
$ws1.Range("F31").Value = 67
$ws1.Range("F32").Value = 329

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14104
$ws4.Range("F5").Value = 550
$ws4.Range("F8").Value = 1049
$ws4.Range("F9").Value = 13934
$ws4.Range("F10").Value = 14985
$ws4.Range("F12").Value = 15
$ws4.Range("F20").Value = 24
$ws4.Range("F21").Value = 63
$ws4.Range("F23").Value = 1157
$ws4.Range("F26").Value = 1
$ws4.Range("F27").Value = 5826
$ws4.Range("F28").Value = 948
$ws4.Range("F29").Value = 1067
$ws4.Range("F30").Value = 5448
$ws4.Range("F33").Value = 67
$ws4.Range("F34").Value = 329
